$d = $word.ActiveDocument

# The date in the first paragraph reads "12/2/2024" and needs to become
# "12/4/2024". The day digit "2" (the single character between the two
# slashes) is the part that changes; it was retyped in place, which is
# why the surrounding text ends up split into three separate runs:
# "12/" + "4" + "/2024".

$dayDigit = $d.Range(3, 4)
$dayDigit.Text = "4"

# Re-select the freshly written character and nudge it with a
# transient bookmark. Adding then immediately removing the bookmark
# forces Word to keep this character as its own run instead of folding
# it back into its identically-formatted neighbours when the document
# is saved.
$dayDigit = $d.Range(3, 4)
$d.Bookmarks.Add("iron_tmp_split", $dayDigit)
$d.Bookmarks("iron_tmp_split").Delete()
